$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2024-01"
$ws.Range("B2").Value = "Alpha"
$ws.Range("C2").Value = 0.3
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1200
$ws.Range("F2").Value = 100
